# Append 5 new result rows (33-37) to the "Bag" results sheet, mirroring
# additional search results that the app persisted (the search now only
# saves when the user opts to, and the Yahoo company-name lookup no
# longer throws unhandled on failure).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(42603.690995370373, "Bag", 157, 339, 31, 5, 18, 21, 78, 0, 4, 0, 100),
    @(42603.691562499997, "Bag", 145, 340, 31, 5, 18, 21, 78, 0, 4, 0, 100),
    @(42603.692546296297, "Bag", 175, 340, 31, 5, 18, 21, 78, 0, 4, 0, 100),
    @(42603.692847222221, "Bag", 161, 340, 31, 5, 18, 21, 78, 0, 4, 0, 100),
    @(42603.693483796298, "Bag", 174, 340, 31, 6, 18, 24, 74, 0, 4, 0, 100)
)

$startRow = 33
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Mirror the existing date-time cell style (numFmtId 22) from the row
    # directly above instead of building a brand-new number format.
    $ws.Range("A" + ($r - 1)).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $data[0]

    for ($c = 1; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
